# Refresh LR-pair data for Pla2g10-Pla2r1 (new TPM values, added Resolving-Mac cluster)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numRows = 16
$numCols = 20
$arr = New-Object 'object[,]' $numRows,$numCols

# Sheet row 2
$rowData = @("ECs", "Pla2g10", "Pla2r1", "ECs", 1, 0.3333333333333333, 0.191414, 0.574242, 0.1152495606707829, 0.1152495606707829, 1, 0.3333333333333333, 0.1060506666666667, 0.318152, 0.007179517892818991, 0.00717951789281899, 0.02029958230933333, 0.182696240784, 0.0008274362829754135, 0.0008274362829754133)
for ($j = 0; $j -lt $numCols; $j++) { $arr[0,$j] = $rowData[$j] }

# Sheet row 3
$rowData = @("ECs", "Pla2g10", "Pla2r1", "FAPs", 1, 0.3333333333333333, 0.191414, 0.574242, 0.1152495606707829, 0.1152495606707829, 3, 1, 7.032353333333333, 21.09706, 0.4760828778567346, 0.4760828778567345, 1.346090880946667, 12.11481792852, 0.05486834251587064, 0.05486834251587063)
for ($j = 0; $j -lt $numCols; $j++) { $arr[1,$j] = $rowData[$j] }

# Sheet row 4
$rowData = @("ECs", "Pla2g10", "Pla2r1", "MuSCs", 1, 0.3333333333333333, 0.191414, 0.574242, 0.1152495606707829, 0.1152495606707829, 3, 1, 7.554679333333333, 22.664038, 0.51144379524419, 0.5114437952441899, 1.446071389910667, 13.014642509196, 0.05894367270969073, 0.05894367270969071)
for ($j = 0; $j -lt $numCols; $j++) { $arr[2,$j] = $rowData[$j] }

# Sheet row 5
$rowData = @("ECs", "Pla2g10", "Pla2r1", "Resolving-Mac", 1, 0.3333333333333333, 0.191414, 0.574242, 0.1152495606707829, 0.1152495606707829, 3, 1, 0.07819633333333333, 0.234589, 0.005293809006256488, 0.005293809006256488, 0.01496787294866667, 0.134710856538, 0.0006101091622460939, 0.0006101091622460938)
for ($j = 0; $j -lt $numCols; $j++) { $arr[3,$j] = $rowData[$j] }

# Sheet row 6
$rowData = @("FAPs", "Pla2g10", "Pla2r1", "ECs", 2, 0.6666666666666666, 0.599873, 1.799619, 0.361180998820695, 0.361180998820695, 1, 0.3333333333333333, 0.1060506666666667, 0.318152, 0.007179517892818991, 0.00717951789281899, 0.06361693156533334, 0.5725523840879999, 0.002593105443579415, 0.002593105443579415)
for ($j = 0; $j -lt $numCols; $j++) { $arr[4,$j] = $rowData[$j] }

# Sheet row 7
$rowData = @("FAPs", "Pla2g10", "Pla2r1", "FAPs", 2, 0.6666666666666666, 0.599873, 1.799619, 0.361180998820695, 0.361180998820695, 3, 1, 7.032353333333333, 21.09706, 0.4760828778567346, 0.4760828778567345, 4.218518891126666, 37.96667002013999, 0.1719520893457264, 0.1719520893457263)
for ($j = 0; $j -lt $numCols; $j++) { $arr[5,$j] = $rowData[$j] }

# Sheet row 8
$rowData = @("FAPs", "Pla2g10", "Pla2r1", "MuSCs", 2, 0.6666666666666666, 0.599873, 1.799619, 0.361180998820695, 0.361180998820695, 3, 1, 7.554679333333333, 22.664038, 0.51144379524419, 0.5114437952441899, 4.531848155724666, 40.786633401522, 0.1847237808069436, 0.1847237808069436)
for ($j = 0; $j -lt $numCols; $j++) { $arr[6,$j] = $rowData[$j] }

# Sheet row 9
$rowData = @("FAPs", "Pla2g10", "Pla2r1", "Resolving-Mac", 2, 0.6666666666666666, 0.599873, 1.799619, 0.361180998820695, 0.361180998820695, 3, 1, 0.07819633333333333, 0.234589, 0.005293809006256488, 0.005293809006256488, 0.04690786906566666, 0.4221708215909999, 0.00191202322444571, 0.001912023224445709)
for ($j = 0; $j -lt $numCols; $j++) { $arr[7,$j] = $rowData[$j] }

# Sheet row 10
$rowData = @("MuSCs", "Pla2g10", "Pla2r1", "ECs", 2, 0.6666666666666666, 0.6933183333333334, 2.079955, 0.417444039211688, 0.417444039211688, 1, 0.3333333333333333, 0.1060506666666667, 0.318152, 0.007179517892818991, 0.00717951789281899, 0.07352687146222223, 0.66174184316, 0.002997046948770947, 0.002997046948770946)
for ($j = 0; $j -lt $numCols; $j++) { $arr[8,$j] = $rowData[$j] }

# Sheet row 11
$rowData = @("MuSCs", "Pla2g10", "Pla2r1", "FAPs", 2, 0.6666666666666666, 0.6933183333333334, 2.079955, 0.417444039211688, 0.417444039211688, 3, 1, 7.032353333333333, 21.09706, 0.4760828778567346, 0.4760828778567345, 4.875659492477777, 43.8809354323, 0.19873795953204, 0.19873795953204)
for ($j = 0; $j -lt $numCols; $j++) { $arr[9,$j] = $rowData[$j] }

# Sheet row 12
$rowData = @("MuSCs", "Pla2g10", "Pla2r1", "MuSCs", 2, 0.6666666666666666, 0.6933183333333334, 2.079955, 0.417444039211688, 0.417444039211688, 3, 1, 7.554679333333333, 22.664038, 0.51144379524419, 0.5114437952441899, 5.237797684254444, 47.14017915829, 0.2134991637164902, 0.2134991637164902)
for ($j = 0; $j -lt $numCols; $j++) { $arr[10,$j] = $rowData[$j] }

# Sheet row 13
$rowData = @("MuSCs", "Pla2g10", "Pla2r1", "Resolving-Mac", 2, 0.6666666666666666, 0.6933183333333334, 2.079955, 0.417444039211688, 0.417444039211688, 3, 1, 0.07819633333333333, 0.234589, 0.005293809006256488, 0.005293809006256488, 0.05421495149944444, 0.487934563495, 0.002209869014386921, 0.002209869014386921)
for ($j = 0; $j -lt $numCols; $j++) { $arr[11,$j] = $rowData[$j] }

# Sheet row 14
$rowData = @("Resolving-Mac", "Pla2g10", "Pla2r1", "ECs", 1, 0.3333333333333333, 0.17626, 0.52878, 0.106125401296834, 0.106125401296834, 1, 0.3333333333333333, 0.1060506666666667, 0.318152, 0.007179517892818991, 0.00717951789281899, 0.01869249050666667, 0.16823241456, 0.0007619292174932156, 0.0007619292174932154)
for ($j = 0; $j -lt $numCols; $j++) { $arr[12,$j] = $rowData[$j] }

# Sheet row 15
$rowData = @("Resolving-Mac", "Pla2g10", "Pla2r1", "FAPs", 1, 0.3333333333333333, 0.17626, 0.52878, 0.106125401296834, 0.106125401296834, 3, 1, 7.032353333333333, 21.09706, 0.4760828778567346, 0.4760828778567345, 1.239522598533333, 11.1557033868, 0.05052448646309757, 0.05052448646309757)
for ($j = 0; $j -lt $numCols; $j++) { $arr[13,$j] = $rowData[$j] }

# Sheet row 16
$rowData = @("Resolving-Mac", "Pla2g10", "Pla2r1", "MuSCs", 1, 0.3333333333333333, 0.17626, 0.52878, 0.106125401296834, 0.106125401296834, 3, 1, 7.554679333333333, 22.664038, 0.51144379524419, 0.5114437952441899, 1.331587779293333, 11.98429001364, 0.05427717801106547, 0.05427717801106547)
for ($j = 0; $j -lt $numCols; $j++) { $arr[14,$j] = $rowData[$j] }

# Sheet row 17
$rowData = @("Resolving-Mac", "Pla2g10", "Pla2r1", "Resolving-Mac", 1, 0.3333333333333333, 0.17626, 0.52878, 0.106125401296834, 0.106125401296834, 3, 1, 0.07819633333333333, 0.234589, 0.005293809006256488, 0.005293809006256488, 0.01378288571333333, 0.12404597142, 0.0005618076051777639, 0.0005618076051777639)
for ($j = 0; $j -lt $numCols; $j++) { $arr[15,$j] = $rowData[$j] }

$ws.Range("A2:T17").Value = $arr
